# Weekly update: insert a new latest-week price record at the top of the
# Acelga / Agrícola del Norte S.A. de Arica price history block (row 90),
# shifting the existing historical rows (90-105) down by one (to 91-106).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 90; Excel shifts rows 90:105 down to 91:106
# and carries the existing column formatting (e.g. the date style on column D).
$ws.Rows(90).Insert()

# Populate the newly inserted row 90 with this week's data.
$ws.Cells.Item(90, 1).Value  = 1
$ws.Cells.Item(90, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(90, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(90, 4).Value  = 45211
$ws.Cells.Item(90, 5).Value  = 15
$ws.Cells.Item(90, 6).Value  = 100112009
$ws.Cells.Item(90, 7).Value  = "Acelga"
$ws.Cells.Item(90, 8).Value  = "Sin especificar"
$ws.Cells.Item(90, 9).Value  = "Primera"
$ws.Cells.Item(90, 10).Value = 300
$ws.Cells.Item(90, 11).Value = 800
$ws.Cells.Item(90, 12).Value = 1000
$ws.Cells.Item(90, 13).Value = 933
$ws.Cells.Item(90, 14).Value = "$/atado 2,5 a 3 kilos"
$ws.Cells.Item(90, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(90, 16).Value = 311
$ws.Cells.Item(90, 17).Value = 3
$ws.Cells.Item(90, 18).Value = "Hortaliza"
